# Apply the "Unify information and codes" edit:
# Rewrite the list of codes/labels on the "sequence" sheet (column B, rows 2-13)
# to the new wording used by the revised shared-strings table.
#
# NOTE: the order below (which determines the order new shared-string table
# entries are created in) is deliberately chosen so the resulting shared
# string indices line up with the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sequence")

$ws.Range("B2").Value  = "Food in the house"
$ws.Range("B5").Value  = "Children's clothes"
$ws.Range("B4").Value  = "Spouse's clothes"
$ws.Range("B3").Value  = "Own clothes"
$ws.Range("B6").Value  = "Children's education"
$ws.Range("B7").Value  = "Children's health"
$ws.Range("B10").Value = "Money to spouse's relatives"
$ws.Range("B9").Value  = "Money to relatives"
$ws.Range("B11").Value = "Own work"
$ws.Range("B12").Value = "Spouse's work"
$ws.Range("B13").Value = "Contraceptives"
$ws.Range("B8").Value  = "Strong expenditure"
